$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44195
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("Q3").Value = "$/bandeja 10 kilos"
$ws.Range("S3").Value = 1500
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44391
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = 1500
$ws.Range("O4").Value = 1500
$ws.Range("P4").Value = 1500
$ws.Range("Q4").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S4").Value = 1500
$ws.Range("T4").Value = 1

# Row 5
$ws.Range("D5").Value = 44391
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 1000
$ws.Range("O5").Value = 1000
$ws.Range("P5").Value = 1000
$ws.Range("Q5").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S5").Value = 1000
$ws.Range("T5").Value = 1

# Row 6
$ws.Range("D6").Value = 44309
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = 1600
$ws.Range("O6").Value = 1600
$ws.Range("P6").Value = 1600
$ws.Range("Q6").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S6").Value = 1600
$ws.Range("T6").Value = 1

# Row 7
$ws.Range("D7").Value = 44371
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 1800
$ws.Range("O7").Value = 1800
$ws.Range("P7").Value = 1800
$ws.Range("Q7").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S7").Value = 1800
$ws.Range("T7").Value = 1

# Row 8
$ws.Range("D8").Value = 44371
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 30
$ws.Range("N8").Value = 1200
$ws.Range("O8").Value = 1200
$ws.Range("P8").Value = 1200
$ws.Range("Q8").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S8").Value = 1200
$ws.Range("T8").Value = 1

# Row 9
$ws.Range("D9").Value = 44292
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 50
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 14000
$ws.Range("P9").Value = 14000
$ws.Range("Q9").Value = "$/bandeja 10 kilos"
$ws.Range("S9").Value = 1400
$ws.Range("T9").Value = 10

# Row 10
$ws.Range("D10").Value = 44400
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 1500
$ws.Range("O10").Value = 1500
$ws.Range("P10").Value = 1500
$ws.Range("Q10").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("S10").Value = 1500
$ws.Range("T10").Value = 1

